# Bewertung.xlsx - update group roster and self-assessment (column D) scores.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Roster update: "Ben Köger" (C5/Mitglied 2) leaves the group and is
# replaced by "Maximilian Marschall". The former "Mitglied 3"
# (Maximilian Marschall) slot becomes "Ben  Köger" (note double space),
# and "Mitglied 4" (previously "Lucas Boyn") is cleared.
$ws.Range("C5").Value = "Maximilian Marschall"
$ws.Range("C6").Value = "Ben  Köger"
$ws.Range("C7").ClearContents()

# Fill in the "Selbsteinschätzung" (self-assessment) column D for each
# grading category.
$ws.Range("D12").Value = 10
$ws.Range("D13").Value = 4
$ws.Range("D14").Value = 9
$ws.Range("D15").Value = 10
$ws.Range("D16").Value = 5
$ws.Range("D17").Value = 5
$ws.Range("D18").Value = 5
$ws.Range("D19").Value = 3

# Leave the active selection on D14, matching the saved view state.
$ws.Range("D14").Select()
